$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content fix: the "Name" column had stray trailing spaces baked into the
# shared strings for two staff members ("Alexei " / "Arvind "). Re-enter the
# values without the trailing space so the cells (and the shared-string
# table) reflect clean names.
$ws.Range("A3").Value = "Alexei"
$ws.Range("A6").Value = "Arvind"

# Selection/view change: the active selection on the frozen (bottom-left)
# pane moved from A7:XFD166 to a single cell, I3.
[void]$ws.Range("I3").Select()

Write-Output "done"
